$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sets")
$ws.Range("A1").Value = "TEST"
